$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Extend existing multi-line cell contents (D4 / E4) with additional items
$ws.Range("D4").Value = "1. CIR-W0010`n2. CIR-W0011`n3. CIR-W0013`n4. CIR-W0014`n5. CIR-W0015`n6. CIR-W0030`n7. CIR-W0031`n8. CIR-W0032`n9. CIR-W0033`n10. CIR-W0034`n11. CIR-W0035"

$ws.Range("E4").Value = "1. Вызов формы захвата с наличием захваченного кадра`n2. Вызов формы захвата без захваченного кадра`n3. Нажатие по кнопке «Старт F5» в форме захвата`n4. Нажатие по кнопке «Стоп F6» в форме захвата`n5. Нажатие по кнопке «Закрыть Esc» в форме захвата`n6. Захват кадра с несколькими лицами в форме захвата`n7. Выбор другого захваченного кадра`n8. Захват кадра на границе области с видеопотоком в форме захвата`n9. Считывание документа с ИС, в открытой форме захвата`n10. Считывание документа без ИС, в открытой форме захвата`n11. Повторное считывание документа с открытой формой захвата во время захвата"

# Fill in new row 8 data (previously empty)
$ws.Range("B8").Value = "CIR-W S2.6 "
$ws.Range("C8").Value = "Проверка логов"
$ws.Range("D8").Value = "1. CIR-W0036`n2. CIR-W0037"
$ws.Range("E8").Value = "1. Создание логов «Контраст» при входе в «Система Каскад»`n2. Изображения в папке «captures»"

# C8 uses the "no border / no wrap" style (same as B3), unlike the other row-8 cells
$ws.Range("C8").WrapText = $false
$ws.Range("C8").Borders.LineStyle = 0

# Adjust row heights
$ws.Rows.Item(4).RowHeight = 329.25
$ws.Rows.Item(5).RowHeight = 47.25
$ws.Rows.Item(8).RowHeight = 47.25

# Update sheet view: top-left cell and selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("E8").Select()
